$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = 25
$ws.Range("B3").Value = 25
$ws.Range("C3").Value = 2000000
$ws.Range("D3").Value = 2000000
$ws.Range("E3").Value = 0.002

$ws.Range("E7").Value = 0.5
$ws.Range("F7").Value = 81.4941

$ws.Range("A12").Value = 3
$ws.Range("B12").Value = 0.000002

$ws.Range("A18").Value = ""
$ws.Range("A19").Value = ""
$ws.Range("A20").Value = ""

$ws.Range("B12").Select() | Out-Null
